# Applies the commit's edits to the document:
#   1. Normalizes the ">>>  your stuff after this line >>>" paragraph -
#      it was split across three runs (with gramStart/gramEnd proofErr
#      markers around "> your") even though the visible text never
#      changed; collapse it back down to a single clean run.
#   2. Splits the "Baz changes" paragraph (itself stored as "Baz chan" +
#      a "_GoBack" bookmark + "ges") into two paragraphs: a clean
#      "Baz changes" line, followed by a new line "Glenn made a change
#      to this file." that now carries the "_GoBack" bookmark at its end
#      (i.e. Glenn's new sentence is the most-recently-edited spot).

$d = $word.ActiveDocument

# --- 1. Re-key the ">>>  your stuff after this line >>>" paragraph -------
# Writing the exact same text back is a no-op for the engine, so first
# push in a throwaway placeholder (forcing a real content change, which
# merges the runs and drops the proofErr markers), then write the real
# text.
$p4 = $d.Paragraphs.Item(4)
$rng = $d.Range($p4.Range.Start, $p4.Range.End - 1)
$rng.Text = "placeholder"
$p4 = $d.Paragraphs.Item(4)
$rng = $d.Range($p4.Range.Start, $p4.Range.End - 1)
$rng.Text = ">>>  your stuff after this line >>>"

# --- 2. Split the "Baz chan" + bookmark + "ges" paragraph ----------------
# Insert a new paragraph break right after it; this paragraph currently
# still reads "Baz changes" (unchanged) and the new paragraph after it
# is empty.
$p5 = $d.Paragraphs.Item(5)
$p5.Range.InsertParagraphAfter()

# Fill the new paragraph with Glenn's sentence, plus one throw-away
# trailing character. Adding a bookmark exactly at the end-of-paragraph
# position is unreliable, so the extra character keeps the bookmark
# insertion point safely *inside* the paragraph; it is deleted afterwards
# and the bookmark (tracked by position) ends up sitting right after the
# real text, which is exactly where it belongs.
$p6 = $d.Paragraphs.Item(6)
$rng = $d.Range($p6.Range.Start, $p6.Range.End - 1)
$rng.Text = "Glenn made a change to this file.X"

$p6 = $d.Paragraphs.Item(6)
$bmPos = $p6.Range.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$p6 = $d.Paragraphs.Item(6)
$trailing = $d.Range($p6.Range.End - 2, $p6.Range.End - 1)
$trailing.Delete()

# --- 3. Clean up the "Baz chan"/"ges" paragraph's run split --------------
# Now that the bookmark lives in the new paragraph, collapse this
# paragraph's remaining two runs back into a single "Baz changes" run the
# same way paragraph 4 was normalized above.
$p5 = $d.Paragraphs.Item(5)
$rng = $d.Range($p5.Range.Start, $p5.Range.End - 1)
$rng.Text = "placeholder2"
$p5 = $d.Paragraphs.Item(5)
$rng = $d.Range($p5.Range.Start, $p5.Range.End - 1)
$rng.Text = "Baz changes"
